$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '51.674.80'
$ws.Range("E2").Value = '  +0.39%  '
$ws.Range("D3").Value = '3.092.06'
$ws.Range("E3").Value = '  +3.63%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '388.84'
$ws.Range("E5").Value = '  +1.77%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '103.33'
$ws.Range("E6").Value = '  -0.80%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("E9").Value = '  -0.65%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '37.07'
$ws.Range("E10").Value = '  +1.12%  '
$ws.Range("E11").Value = '  +0.16%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0861'
$ws.Range("E12").Value = '  +0.38%  '
$ws.Range("D13").Value = '3.579.34'
$ws.Range("E13").Value = '  +3.56%  '
$ws.Range("E14").Value = '  +1.23%  '
$ws.Range("E15").Value = '  -0.03%  '
$ws.Range("D16").Value = '3.089.58'
$ws.Range("E16").Value = '  +3.52%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.983'
$ws.Range("E17").Value = '  -1.29%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '10.66'
$ws.Range("E18").Value = '  -4.28%  '
$ws.Range("D19").Value = '51.812.66'
$ws.Range("E19").Value = '  +0.63%  '
$ws.Range("E20").Value = '  +3.03%  '
$ws.Range("E21").Value = '  -0.85%  '
$ws.Range("E22").Value = '  +0.82%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '70.03'
$ws.Range("E23").Value = '  -0.34%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '268.80'
$ws.Range("E24").Value = '  +0.66%  '
$ws.Range("E25").Value = '  -2.47%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.18'
$ws.Range("E26").Value = '  +4.13%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '27.10'
$ws.Range("E27").Value = '  +4.02%  '
$ws.Range("E28").Value = '  +0.42%  '
$ws.Range("E29").Value = '  -1.76%  '
$ws.Range("E30").Value = '  +0.05%  '
$ws.Range("E31").Value = '  -0.44%  '
$ws.Range("E32").Value = '  -0.55%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '35.63'
$ws.Range("E33").Value = '  +2.83%  '
$ws.Range("E34").Value = '  -1.96%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '50.48'
$ws.Range("E36").Value = '  +0.30%  '
$ws.Range("E37").Value = '  -0.11%  '
$ws.Range("E38").Value = '  +3.45%  '
$ws.Range("E39").Value = '  +8.19%  '
$ws.Range("B40").Value = 'ARBITRUM'
$ws.Range("C40").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.88'
$ws.Range("E40").Value = '  +2.23%  '
$ws.Range("B41").Value = 'Celestia'
$ws.Range("C41").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '17.02'
$ws.Range("E41").Value = '  +0.35%  '
$ws.Range("E42").Value = '  +0.47%  '
$ws.Range("E43").Value = '  -0.35%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '127.34'
$ws.Range("E44").Value = '  +0.26%  '
$ws.Range("E45").Value = '  -2.67%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '22.14'
$ws.Range("E46").Value = '  +3.43%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.49'
$ws.Range("E47").Value = '  +5.19%  '
$ws.Range("E48").Value = '  +2.47%  '
$ws.Range("D49").Value = '2.046.50'
$ws.Range("E49").Value = '  +1.11%  '
$ws.Range("D50").Value = '3.397.38'
$ws.Range("E50").Value = '  +3.64%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.207'
$ws.Range("E51").Value = '  +6.77%  '
